# Rename the data sheet used by the Cucumber test-data workbook from the
# default "Sheet1" to "Credentials" so the feature file / step-defs can
# reference it by its intended name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Credentials"
